$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "FT232500M4638RWZ"
